# Updated cryptos list on Tue Jun 18 23:23:43 UTC 2024 with GitHub Actions
# Refresh the Price (col D) / Volume(1h) (col E) columns with the latest
# scraped figures, and re-sort a handful of near-tied rows (Coin/Link/
# Price/Volume all move together for those).
#
# Cells D/E are stored as text (not numbers) in the sheet, so for any new
# value that Excel would otherwise auto-detect as numeric we briefly force
# the cell to Text format, assign the value, then restore the Normal style
# so formatting matches the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "65.158.74"
$ws.Cells.Item(2, 5).Value = "  -2.04%  "
$ws.Cells.Item(3, 4).Value = "3.478.97"
$ws.Cells.Item(3, 5).Value = "  -0.92%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.11%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "588.98"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -2.76%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "137.24"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -4.13%  "
$ws.Cells.Item(7, 4).Value = "3.478.57"
$ws.Cells.Item(7, 5).Value = "  -0.87%  "
$ws.Cells.Item(8, 5).Value = "  +0.08%  "
$ws.Cells.Item(9, 5).Value = "  -2.86%  "
$ws.Cells.Item(10, 5).Value = "  -5.50%  "
$ws.Cells.Item(11, 5).Value = "  -6.82%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.384"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -4.67%  "
$ws.Cells.Item(13, 4).Value = "4.066.57"
$ws.Cells.Item(13, 5).Value = "  -0.82%  "
$ws.Cells.Item(14, 5).Value = "  -6.31%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "26.63"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -6.83%  "
$ws.Cells.Item(16, 4).Value = "3.451.75"
$ws.Cells.Item(16, 5).Value = "  -1.46%  "
$ws.Cells.Item(17, 5).Value = "  -1.32%  "
$ws.Cells.Item(18, 4).Value = "65.092.07"
$ws.Cells.Item(18, 5).Value = "  -1.96%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "9.73"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -8.48%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "5.78"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -5.13%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "13.92"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -4.42%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "389.41"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -7.22%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.556"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -5.05%  "
$ws.Cells.Item(24, 2).Value = "Litecoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "72.60"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -5.39%  "
$ws.Cells.Item(25, 2).Value = "Dai"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.00"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.04%  "
$ws.Cells.Item(26, 2).Value = "WrappedeETH"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(26, 4).Value = "3.621.24"
$ws.Cells.Item(26, 5).Value = "  -0.98%  "
$ws.Cells.Item(27, 2).Value = "LEO"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "5.75"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.27%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.0000110"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -2.05%  "
$ws.Cells.Item(29, 5).Value = "  +0.11%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "7.35"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -5.13%  "
$ws.Cells.Item(31, 5).Value = "  -8.15%  "
$ws.Cells.Item(32, 5).Value = "  -9.55%  "
$ws.Cells.Item(33, 4).Value = "3.497.92"
$ws.Cells.Item(33, 5).Value = "  -0.54%  "
$ws.Cells.Item(34, 5).Value = "  -0.01%  "
$ws.Cells.Item(35, 2).Value = "Kaspa"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.143"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -6.68%  "
$ws.Cells.Item(36, 2).Value = "EthereumClassic"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "23.11"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -4.21%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "171.22"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.98%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.20"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -8.54%  "
$ws.Cells.Item(39, 5).Value = "  -9.00%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "1.48"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -8.63%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "4.72"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -8.52%  "
$ws.Cells.Item(42, 5).Value = "  -2.70%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.812"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -4.48%  "
$ws.Cells.Item(44, 2).Value = "OKB"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "42.50"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -6.59%  "
$ws.Cells.Item(45, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.00"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.10%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "25.24"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +11.24%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "4.36"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -11.97%  "
$ws.Cells.Item(48, 5).Value = "  -7.74%  "
$ws.Cells.Item(49, 5).Value = "  +4.06%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "6.70"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -4.61%  "
$ws.Cells.Item(51, 2).Value = "dogwifhat"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.06"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -12.07%  "
